# Discharge workbook: add cross-sectional Area / Atotal columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2/C2 used to hold the placeholder string "-" (no velocity measurement at
# the bank, x=30). They become real numeric zeros so the new Area formulas
# (which multiply B by a width) evaluate cleanly instead of propagating text.
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0

# New headers.
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"

# Per-segment cross-sectional area.
# Row 2 (left bank, x=30) anchors against 0 instead of the previous station.
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
# Row 3 is its own one-off (D3-D2) before the fill-down pattern stabilizes.
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"
# Rows 4-15 share the same relative formula -> Excel stores this as one
# shared-formula group, same as a fill-down/fill-to would.
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

# Running total of area across the cross-section.
$ws.Range("H2").Formula = "=SUM(G2:G11)"

# Restore the selection to match the freshly-edited cell.
$ws.Range("H2").Select() | Out-Null
